$d = $word.ActiveDocument

# 1. "From author annotations" -> "For author annotations"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*From author annotations*") {
        $p.Range.Text = "For author annotations"
    }
}

# 2. Add two new list items ("Credit Title" and "Credit Authors") right after
#    the "Resource Title" bullet, matching its list level/numbering.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Resource Title*") {
        $resourceTitle = $p
    }
}
$resourceTitle.Range.InsertAfter("`rCredit Title`rCredit Authors")
